# AE-5 Estudar o modelo de dados definido para a RSL
# Refresh the sprint metrics table: update rows 2-6 with new values and
# append new sprint rows 7-11 (Sprint_06 .. Sprint_10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2 = @("Sprint_01"; 45658; 45672; 5; 14; 41; 40; 79; 64; "Web"; 4.1; 0.17; 0.08; 66.98; 15.8)
    3 = @("Sprint_02"; 45672; 45686; 3; 14; 41; 39; 76; 60; "Dados"; 2.6; 0.1; 0.13; 34.02; 25.33)
    4 = @("Sprint_03"; 45686; 45700; 7; 14; 23; 22; 86; 85; "Mobile"; 3.2; 0.19; 0.16; 68.04; 12.29)
    5 = @("Sprint_04"; 45700; 45714; 6; 14; 46; 40; 75; 72; "Dados"; 4.1; 0.19; 0.16; 56.08; 12.5)
    6 = @("Sprint_05"; 45714; 45728; 7; 14; 34; 33; 92; 91; "Dados"; 2.7; 0.15; 0.16; 33.97; 13.14)
    7 = @("Sprint_06"; 45728; 45742; 7; 14; 44; 34; 109; 104; "Dados"; 3; 0.07; 0.06; 56.34; 15.57)
    8 = @("Sprint_07"; 45742; 45756; 7; 14; 29; 28; 104; 102; "Mobile"; 2.8; 0.23; 0.1; 45.68; 14.86)
    9 = @("Sprint_08"; 45756; 45770; 6; 14; 45; 41; 83; 77; "Dados"; 1.6; 0.23; 0.14; 62.75; 13.83)
    10 = @("Sprint_09"; 45770; 45784; 5; 14; 28; 23; 95; 70; "API"; 2.1; 0.17; 0.16; 64.04; 19)
    11 = @("Sprint_10"; 45784; 45798; 4; 14; 50; 46; 65; 53; "Dados"; 4.5; 0.18; 0.17; 48.75; 16.25)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item($r, $col).Value = $vals[$i]
    }
}

# Apply the date number format (yyyy-mm-dd, matching the existing B/C date
# columns) to the newly appended rows.
$ws.Range("B7:C11").NumberFormat = "YYYY-MM-DD"
